$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" row query (B2) is corrected: the trailing
# "coalesce(co.cohort_description, '') AS `Cohort`" return column is removed
# (it referenced a variable/field that wasn't actually needed/valid for this query).
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN['Golden Retriever']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Reset the sheet selection/view back to the top of the sheet with B2 as the
# active cell (instead of being scrolled down to B4).
$ws.Range("B2").Select()

$wb.Save()
